$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '41.232.02'
$ws.Range('E2').Value = '  -1.89%  '
$ws.Range('D3').Value = '2.174.89'
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '237.06'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -2.26%  '
$ws.Range('E6').Value = '  -2.14%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '70.32'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -5.12%  '
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('E9').Value = '  -5.99%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '40.16'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -9.39%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0929'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -3.43%  '
$ws.Range('E12').Value = '  -5.30%  '
$ws.Range('E13').Value = '  -2.20%  '
$ws.Range('D14').Value = '2.500.46'
$ws.Range('E14').Value = '  -2.05%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '13.94'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -2.52%  '
$ws.Range('E16').Value = '  -4.51%  '
$ws.Range('D17').Value = '2.170.56'
$ws.Range('E17').Value = '  -2.20%  '
$ws.Range('D18').Value = '41.017.97'
$ws.Range('E18').Value = '  -2.15%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.0000102'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -7.14%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '70.44'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -2.81%  '
$ws.Range('E21').Value = '  -3.94%  '
$ws.Range('E22').Value = '  -10.99%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '225.98'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.94%  '
$ws.Range('E24').Value = '  -7.47%  '
$ws.Range('E25').Value = '  +0.13%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '10.83'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -6.61%  '
$ws.Range('E27').Value = '  -0.82%  '
$ws.Range('E28').Value = '  -2.89%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.19'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.69%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '167.25'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.35%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '19.99'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -3.21%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '31.02'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +6.58%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0771'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -3.89%  '
$ws.Range('E34').Value = '  -9.23%  '
$ws.Range('E35').Value = '  -3.11%  '
$ws.Range('E36').Value = '  -8.64%  '
$ws.Range('E37').Value = '  -4.40%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0286'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -5.66%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '12.17'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E40').Value = '  -2.92%  '
$ws.Range('E41').Value = '  -4.01%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '60.48'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -8.15%  '
$ws.Range('E43').Value = '  -4.80%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0973'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -3.77%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '98.23'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -5.89%  '
$ws.Range('E47').Value = '  -2.80%  '
$ws.Range('E48').Value = '  -2.55%  '
$ws.Range('E49').Value = '  -8.52%  '
$ws.Range('E50').Value = '  -3.05%  '
$ws.Range('D51').Value = '2.379.10'
$ws.Range('E51').Value = '  -2.01%  '
